$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 1010
$ws.Cells.Item(7, 6).Value = 2537
$ws.Cells.Item(9, 6).Value = 1243
$ws.Cells.Item(10, 6).Value = 907
$ws.Cells.Item(11, 6).Value = 604
$ws.Cells.Item(12, 6).Value = 912
$ws.Cells.Item(13, 6).Value = 1130
$ws.Cells.Item(16, 6).Value = 111
$ws.Cells.Item(18, 6).Value = 772
$ws.Cells.Item(19, 6).Value = 193
$ws.Cells.Item(20, 6).Value = 482
$ws.Cells.Item(21, 6).Value = 1111
$ws.Cells.Item(23, 6).Value = 599
$ws.Cells.Item(24, 6).Value = 589
$ws.Cells.Item(29, 6).Value = 392
$ws.Cells.Item(30, 6).Value = 4329
$ws.Cells.Item(31, 6).Value = 481
$ws.Cells.Item(36, 6).Value = 146
$ws.Cells.Item(37, 6).Value = 1601
$ws.Cells.Item(40, 6).Value = 84
$ws.Cells.Item(41, 6).Value = 141
$ws.Cells.Item(46, 6).Value = 100
$ws.Cells.Item(47, 6).Value = 22

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(6, 6).Value = 30
$ws.Cells.Item(8, 6).Value = 18
$ws.Cells.Item(13, 6).Value = 20
$ws.Cells.Item(16, 6).Value = 185

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 2261
$ws.Cells.Item(3, 6).Value = 727

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 2261
$ws.Cells.Item(3, 6).Value = 727
$ws.Cells.Item(7, 6).Value = 1010
$ws.Cells.Item(8, 6).Value = 2537
$ws.Cells.Item(10, 6).Value = 1243
$ws.Cells.Item(11, 6).Value = 907
$ws.Cells.Item(12, 6).Value = 604
$ws.Cells.Item(13, 6).Value = 912
$ws.Cells.Item(14, 6).Value = 1130
$ws.Cells.Item(17, 6).Value = 111
$ws.Cells.Item(19, 6).Value = 772
$ws.Cells.Item(20, 6).Value = 193
$ws.Cells.Item(21, 6).Value = 482
$ws.Cells.Item(22, 6).Value = 1112
$ws.Cells.Item(25, 6).Value = 599
$ws.Cells.Item(26, 6).Value = 589
$ws.Cells.Item(30, 6).Value = 393
$ws.Cells.Item(31, 6).Value = 4329
$ws.Cells.Item(33, 6).Value = 481
$ws.Cells.Item(36, 6).Value = 146
$ws.Cells.Item(37, 6).Value = 1601
$ws.Cells.Item(39, 6).Value = 20
$ws.Cells.Item(40, 6).Value = 20
$ws.Cells.Item(42, 6).Value = 84
$ws.Cells.Item(43, 6).Value = 141
$ws.Cells.Item(48, 6).Value = 100
